$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 15; existing rows 15..43 shift down to 16..44.
$ws.Rows.Item(15).Insert()

# Populate the new row 15 with the new weekly record (Berenjena, Vega Monumental
# Concepcion). Text columns as strings, numeric columns as numbers; column D
# keeps the date-formatted style (already inherited from the row insert).
$ws.Range("A15").Value = 11
$ws.Range("B15").Value = "Vega Monumental Concepción"
$ws.Range("C15").Value = "Bíobío"
$ws.Range("D15").Value = 44497
$ws.Range("E15").Value = 8
$ws.Range("F15").Value = 100112001
$ws.Range("G15").Value = "Berenjena"
$ws.Range("H15").Value = "Sin especificar"
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 220
$ws.Range("K15").Value = 7500
$ws.Range("L15").Value = 8000
$ws.Range("M15").Value = 7727
$ws.Range("N15").Value = "$/caja 60 unidades"
$ws.Range("O15").Value = "Región de Arica y Parinacota"
$ws.Range("P15").Value = 129
$ws.Range("Q15").Value = 60
$ws.Range("R15").Value = "Hortaliza"
